$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "67.319.74"
$ws.Cells.Item(2, 5).Value = "  +7.49%  "
$ws.Cells.Item(3, 4).Value = "3.589.90"
$ws.Cells.Item(3, 5).Value = "  +3.72%  "
$ws.Cells.Item(4, 5).Value = "  +0.00%  "
$ws.Cells.Item(5, 4).Value = "'417.14"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.83%  "
$ws.Cells.Item(6, 4).Value = "'129.30"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.42%  "
$ws.Cells.Item(7, 4).Value = "'0.653"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +4.31%  "
$ws.Cells.Item(8, 4).Value = "3.580.85"
$ws.Cells.Item(8, 5).Value = "  +3.59%  "
$ws.Cells.Item(9, 5).Value = "  -0.01%  "
$ws.Cells.Item(10, 4).Value = "'0.776"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +6.92%  "
$ws.Cells.Item(11, 5).Value = "  +16.17%  "
$ws.Cells.Item(12, 5).Value = "  +51.49%  "
$ws.Cells.Item(13, 4).Value = "'42.60"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +0.21%  "
$ws.Cells.Item(14, 4).Value = "'9.91"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +2.20%  "
$ws.Cells.Item(15, 4).Value = "4.165.75"
$ws.Cells.Item(15, 5).Value = "  +3.80%  "
$ws.Cells.Item(16, 5).Value = "  -0.12%  "
$ws.Cells.Item(17, 4).Value = "'20.45"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -0.58%  "
$ws.Cells.Item(18, 4).Value = "3.610.99"
$ws.Cells.Item(18, 5).Value = "  +3.68%  "
$ws.Cells.Item(19, 5).Value = "  +6.02%  "
$ws.Cells.Item(20, 4).Value = "67.206.00"
$ws.Cells.Item(20, 5).Value = "  +7.36%  "
$ws.Cells.Item(21, 4).Value = "'12.31"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -3.13%  "
$ws.Cells.Item(22, 4).Value = "'452.89"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -2.19%  "
$ws.Cells.Item(23, 4).Value = "'89.59"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -1.04%  "
$ws.Cells.Item(24, 4).Value = "'3.16"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -3.39%  "
$ws.Cells.Item(25, 4).Value = "'13.19"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.69%  "
$ws.Cells.Item(26, 5).Value = "  +1.64%  "
$ws.Cells.Item(27, 4).Value = "'10.04"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -5.09%  "
$ws.Cells.Item(28, 4).Value = "'35.06"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +5.28%  "
$ws.Cells.Item(29, 4).Value = "'4.88"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +1.94%  "
$ws.Cells.Item(30, 2).Value = "Toncoin"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(30, 4).Value = "'2.78"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +4.19%  "
$ws.Cells.Item(31, 2).Value = "Cosmos"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(31, 4).Value = "'12.42"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +3.19%  "
$ws.Cells.Item(32, 5).Value = "  +5.03%  "
$ws.Cells.Item(33, 5).Value = "  -2.02%  "
$ws.Cells.Item(34, 4).Value = "'0.163"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -2.25%  "
$ws.Cells.Item(35, 4).Value = "'41.62"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +1.98%  "
$ws.Cells.Item(36, 4).Value = "'0.999"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +0.03%  "
$ws.Cells.Item(37, 4).Value = "'56.67"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -3.03%  "
$ws.Cells.Item(38, 4).Value = "'0.0495"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +0.91%  "
$ws.Cells.Item(39, 4).Value = "0.0₃0752"
$ws.Cells.Item(39, 5).Value = "  +33.96%  "
$ws.Cells.Item(40, 4).Value = "'0.147"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +9.84%  "
$ws.Cells.Item(41, 2).Value = "Stacks"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(41, 4).Value = "'3.09"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.38%  "
$ws.Cells.Item(42, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(42, 4).Value = "'0.998"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.22%  "
$ws.Cells.Item(43, 2).Value = "WEMIXToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(43, 4).Value = "'2.75"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +1.48%  "
$ws.Cells.Item(44, 2).Value = "Monero"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(44, 4).Value = "'148.92"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.47%  "
$ws.Cells.Item(45, 4).Value = "'0.318"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -1.05%  "
$ws.Cells.Item(46, 4).Value = "'3.28"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -1.47%  "
$ws.Cells.Item(47, 4).Value = "'4.34"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -1.37%  "
$ws.Cells.Item(48, 5).Value = "  -4.10%  "
$ws.Cells.Item(49, 4).Value = "'2.32"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -3.45%  "
$ws.Cells.Item(50, 2).Value = "Celestia"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(50, 4).Value = "'15.73"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -4.16%  "
$ws.Cells.Item(51, 2).Value = "BitcoinSV"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Cells.Item(51, 4).Value = "'115.56"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +6.56%  "
